$d = $word.ActiveDocument

# Locate the run containing "Shoaib Ahmed (F17040117" so we can fix the
# typo'd roll number and relocate the "_GoBack" bookmark there (Word
# keeps "_GoBack" pinned to the most recent edit location).
$rng = $d.Content
$found = $rng.Find.Execute("F17040117")
if (-not $found) {
    throw "Could not find target text 'F17040117'"
}

$matchEnd = $rng.End

# Re-seat the "_GoBack" bookmark right after the matched run (this also
# removes it from its old spot after "Islam", since bookmark names are
# unique per document). Doing this before the text edit keeps the run
# split exactly where Word would leave it (edit point == bookmark point).
$d.Bookmarks.Add("_GoBack", $d.Range($matchEnd, $matchEnd))

# Fix the last digit of the roll number: 117 -> 116.
$lastChar = $d.Range($matchEnd - 1, $matchEnd)
$lastChar.Text = "6"
